$d = $word.ActiveDocument
Write-Output "test"
